$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.303.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.930.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.72%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7264"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3336"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.60%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "28.01"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06921"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8058"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08041"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.933.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.420"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.87%  "

$ws.Range("E16").Value = "  -2.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.304.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "253.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -8.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008200"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.806"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.187.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9993"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.888"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.748"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.456"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1335"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.88%  "

$ws.Range("E30").Value = "  -4.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.339"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.418"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.201"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05148"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.225"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7449"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.744"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01991"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.832"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.643"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "78.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4476"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.002"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8372"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.789"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.322"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.489"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.80%  "

$ws.Range("E51").Value = "  -0.46%  "
